$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old C1 cell content (column C no longer used)
$ws.Range("C1").ClearContents()

# Update row 1
$ws.Range("A1").Value = "under pressure"
$ws.Range("B1").Value = "david bowie"

# Add row 2
$ws.Range("A2").Value = "norwegian wood"
$ws.Range("B2").Value = "the beatles"

# Set column A width to match bestFit (13.6640625) via AutoFit
$ws.Columns.Item(1).AutoFit() | Out-Null

# Update selection to A3 as shown in the diff
$ws.Range("A3").Select()
